$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.462.02'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.573.17'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  -0.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.98'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3408'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07567'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.142'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.45%  '
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.991'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.944'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = '1.569.76'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001123'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06733'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.273'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.19'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.39%  '
$ws.Range('D24').Value = '22.466.41'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.337'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.582'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.40'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.011'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.93'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Value = '1.746.17'
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.046'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.121'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.979'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.835'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08400'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.383'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02458'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2297'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06523'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.481'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.36'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6277'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.97'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.812'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5856'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.087'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '129.68'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.231'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07332'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.05%  '
